# Array Problems Day 2 commit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Array")

# Mark "Completed" (column D) for rows 3, 9 and 11 with "Y"
$ws.Range("D3").Value = "Y"
$ws.Range("D9").Value = "Y"
$ws.Range("D11").Value = "Y"

# Update the selection to D13 (single cell)
$ws.Activate()
$ws.Range("D13").Select()
